$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '26.037.39'
$ws.Cells.Item(2, 5).Value = '  +0.58%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '1.644.41'
$ws.Cells.Item(3, 5).Value = '  +0.84%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.70%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '216.40'
$ws.Cells.Item(5, 5).Value = '  +0.80%  '

# Row 6
$ws.Cells.Item(6, 5).Value = '  +0.99%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  +0.67%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  +0.56%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  +1.21%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '19.67'
$ws.Cells.Item(10, 5).Value = '  +0.10%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.0796'
$ws.Cells.Item(11, 5).Value = '  +1.25%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '1.872.99'
$ws.Cells.Item(12, 5).Value = '  +0.89%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '1.659.30'
$ws.Cells.Item(14, 5).Value = '  +4.74%  '

# Row 15
$ws.Cells.Item(15, 5).Value = '  +0.06%  '

# Row 16
$ws.Cells.Item(16, 5).Value = '  +0.97%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '63.24'
$ws.Cells.Item(17, 5).Value = '  +0.72%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '26.046.11'
$ws.Cells.Item(18, 5).Value = '  +0.64%  '

# Row 19
$ws.Cells.Item(19, 5).Value = '  +0.69%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '193.14'
$ws.Cells.Item(20, 5).Value = '  +0.25%  '

# Row 21
$ws.Cells.Item(21, 5).Value = '  -0.66%  '

# Row 22
$ws.Cells.Item(22, 5).Value = '  +0.07%  '

# Row 23
$ws.Cells.Item(23, 5).Value = '  -0.12%  '

# Row 24
$ws.Cells.Item(24, 5).Value = '  +5.13%  '

# Row 25
$ws.Cells.Item(25, 2).Value = 'Monero'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '144.58'
$ws.Cells.Item(25, 5).Value = '  +1.48%  '

# Row 26
$ws.Cells.Item(26, 2).Value = 'Toncoin'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '1.80'
$ws.Cells.Item(26, 5).Value = '  +0.66%  '

# Row 27
$ws.Cells.Item(27, 5).Value = '  +0.82%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '6.90'
$ws.Cells.Item(28, 5).Value = '  +0.59%  '

# Row 29
$ws.Cells.Item(29, 5).Value = '  +0.74%  '

# Row 30
$ws.Cells.Item(30, 5).Value = '  +1.07%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '0.0499'
$ws.Cells.Item(31, 5).Value = '  +0.22%  '

# Row 32
$ws.Cells.Item(32, 5).Value = '  -0.52%  '

# Row 33
$ws.Cells.Item(33, 5).Value = '  +1.20%  '

# Row 34
$ws.Cells.Item(34, 5).Value = '  -2.87%  '

# Row 35
$ws.Cells.Item(35, 5).Value = '  +2.46%  '

# Row 36
$ws.Cells.Item(36, 5).Value = '  +0.61%  '

# Row 37
$ws.Cells.Item(37, 4).Value = '1.132.13'
$ws.Cells.Item(37, 5).Value = '  -0.30%  '

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.541'
$ws.Cells.Item(38, 5).Value = '  -1.52%  '

# Row 39
$ws.Cells.Item(39, 5).Value = '  +0.37%  '

# Row 40
$ws.Cells.Item(40, 5).Value = '  +0.75%  '

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '5.51'
$ws.Cells.Item(41, 5).Value = '  +0.74%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '99.49'
$ws.Cells.Item(42, 5).Value = '  +0.25%  '

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.798'
$ws.Cells.Item(43, 5).Value = '  -0.54%  '

# Row 44
$ws.Cells.Item(44, 4).Value = '1.782.55'
$ws.Cells.Item(44, 5).Value = '  +0.94%  '

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '56.56'
$ws.Cells.Item(46, 5).Value = '  +0.96%  '

# Row 47
$ws.Cells.Item(47, 5).Value = '  -0.44%  '

# Row 48
$ws.Cells.Item(48, 5).Value = '  +0.19%  '

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '7.70'
$ws.Cells.Item(49, 5).Value = '  +1.07%  '

# Row 50
$ws.Cells.Item(50, 5).Value = '  +0.41%  '

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.0957'
$ws.Cells.Item(51, 5).Value = '  -0.42%  '
